# Week 15 simulations update
# OFF sheet ("H" row, row 2): Short Att/Comp, Deep Att/Comp, Short Int, Deep Int
$wsOff = $excel.ActiveWorkbook.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 351
$wsOff.Range("C2").Value = 239
$wsOff.Range("D2").Value = 78
$wsOff.Range("E2").Value = 34
$wsOff.Range("G2").Value = 6

# DEF sheet ("H" row, row 2)
$wsDef = $excel.ActiveWorkbook.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 516
$wsDef.Range("C2").Value = 370
$wsDef.Range("D2").Value = 137
$wsDef.Range("E2").Value = 67
